# [unc-ebike] add spanish translations
#
# Updates the Spanish ("label::Español (es)") translations on the "survey"
# and "choices" sheets of the onboarding XLSForm, and leaves the workbook
# selection/active-sheet state the way the authoring session left it.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("survey")
$ws2 = $wb.Worksheets.Item("choices")

# --- Update Spanish translations -----------------------------------------
# Order matters here only insofar as it determines the order new shared
# strings get appended in; set the "choices" sheet translations first so
# they match the authoring order, then the "survey" sheet question label.

# choices!D4 - "dontknow" Spanish label: "No lo sé" -> "No sé"
$ws2.Range("D4").Value = "No sé"

# choices!D5 - "other" Spanish label: "Otro tipo" -> "Otro"
$ws2.Range("D5").Value = "Otro"

# survey!D2 - smartphone question Spanish label, re-worded
$ws1.Range("D2").Value = "¿Qué tipo de teléfono inteligente usas?"

# The re-worded question label is longer, so wrap the text and let the row
# grow to fit it, same as Excel does when you turn on wrap text for a
# multi-line label.
$ws1.Range("D2").WrapText = $true
$ws1.Rows.Item(2).RowHeight = 68

# --- Restore the selection / active sheet state left by the edit ---------
# The original file had "choices" as the active tab with C8 selected; the
# saved file instead has "survey" active with H2 selected, and "choices"
# showing F6 selected in the background.
[void]$ws2.Range("F6").Select()
[void]$ws1.Activate()
[void]$ws1.Range("H2").Select()
